# Weekly fruit/vegetable price update: a new week's record is inserted at
# row 13 (pushing the two existing rows down to 14 and 15).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13; existing rows 13 and 14 shift to 14 and 15.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with this week's data.
$ws.Cells.Item(13, 1).Value = 1
$ws.Cells.Item(13, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(13, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(13, 4).Value = 44524
$ws.Cells.Item(13, 5).Value = 15
$ws.Cells.Item(13, 6).Value = 100112003
$ws.Cells.Item(13, 7).Value = "Ajo"
$ws.Cells.Item(13, 8).Value = "Chino"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 200
$ws.Cells.Item(13, 11).Value = 20000
$ws.Cells.Item(13, 12).Value = 21000
$ws.Cells.Item(13, 13).Value = 20500
$ws.Cells.Item(13, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(13, 15).Value = "China"
$ws.Cells.Item(13, 16).Value = 2050
$ws.Cells.Item(13, 17).Value = 10
$ws.Cells.Item(13, 18).Value = "Hortaliza"
